# Fixed empty string issue
#
# Sheet1: a new "Time Unit" column is inserted before column D (shifting the
# old Axis/Title/Range columns one to the right), and the row that used to
# hold a bogus/empty-string "e" marker is corrected to a real "D" input with
# its Format/Axis/Range values filled in.
#
# Sheet2: the "Time Unit" / "Time Axis" columns (H:I) are removed entirely.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert the new "Time Unit" column before column D; this shifts the
# existing D/E/F columns (Axis/Title/Range) one column to the right.
$ws1.Range("D:D").EntireColumn.Insert()

$ws1.Range("D1").Value = "Time Unit"

# The old row 4 referenced an empty/placeholder "e" string; fix it up to a
# proper "D" input row with its Format + Axis values filled in.
$ws1.Range("A4").Value = "D"
$ws1.Range("C4").Value = 1
$ws1.Range("E4").Value = "y"

# That "y" axis value used to live (incorrectly) on row 5 - remove it there.
$ws1.Range("E5").ClearContents()

# Row 6 gains an explicit "D" Time Unit value.
$ws1.Range("D6").Value = "D"

# The stray "4:10" range is corrected to "5:20" (matching rows 5 & 6).
$ws1.Range("G4").Value = "5:20"

$ws1.Range("A1:G1").Select()

$ws2 = $wb.Worksheets.Item("Sheet2")

# Drop the "Time Unit" / "Time Axis" columns entirely.
$ws2.Range("H:I").EntireColumn.Delete()

$ws2.Range("H1:J2").Select()
